# "aggiornamento fino a 28 luglio" - append new daily rows (302-328) to Sheet1.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, date-serial (col A), nuovi pos. (col B), somma mobile 7gg. (col C), somma mobile 7gg. per 100mila abitanti (col D)
$newData = @(
    @(302, 44376, 0, 0, 0),
    @(303, 44377, 0, 0, 0),
    @(304, 44378, 0, 0, 0),
    @(305, 44379, 0, 0, 0),
    @(306, 44380, 0, 0, 0),
    @(307, 44381, 0, 0, 0),
    @(308, 44382, 0, 0, 0),
    @(309, 44383, 0, 0, 0),
    @(310, 44384, 0, 0, 0),
    @(311, 44385, 0, 0, 0),
    @(312, 44386, 0, 0, 0),
    @(313, 44387, 0, 0, 0),
    @(314, 44388, 0, 0, 0),
    @(315, 44389, 0, 0, 0),
    @(316, 44390, 0, 0, 0),
    @(317, 44391, 0, 0, 0),
    @(318, 44392, 0, 0, 0),
    @(319, 44393, 0, 0, 0),
    @(320, 44394, 0, 0, 0),
    @(321, 44395, 1, 1, 16.02307322544464),
    @(322, 44396, 1, 2, 32.04614645088928),
    @(323, 44397, 0, 2, 32.04614645088928),
    @(324, 44398, 0, 2, 32.04614645088928),
    @(325, 44399, 0, 2, 32.04614645088928),
    @(326, 44400, 0, 2, 32.04614645088928),
    @(327, 44401, 0, 2, 32.04614645088928),
    @(328, 44402, 1, 2, 32.04614645088928)
)

foreach ($entry in $newData) {
    $r = $entry[0]
    $ws.Cells.Item($r, 1).Value = $entry[1]
    $ws.Cells.Item($r, 2).Value = $entry[2]
    $ws.Cells.Item($r, 3).Value = $entry[3]
    $ws.Cells.Item($r, 4).Value = $entry[4]
}

# Apply the same formatting (bold, centered, bordered, date number format) used by
# the existing date column (A2:A301) to the newly added date cells.
$ws.Range("A301").Copy() | Out-Null
$ws.Range("A302:A328").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
